$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "Bom Dia Inter"
$ws.Cells.Item($row, 3).Value = "Infraestrutura"
$ws.Cells.Item($row, 4).Value = "2025-04-04T08:23"
$ws.Cells.Item($row, 5).Value = "Negativo"
$ws.Cells.Item($row, 6).Value = "Desabou e ficou por isso mesmo. Mureta do Canal do Saco cedeu há meses no Parque Rodoviário, em Campos. Repórter *ao vivo* do local. Situação preocupa moradores e motoristas. Parte da estrutura desabou em janeiro. Até hoje, nada foi feito. Sinal interrompido. Depoimento de morador próximo. Trecho movimentado porque tem creche e escola ali perto. Equipe entrou em contato com a prefeitura e aguarda posicionamento. "
